# Append a new log row (row 10) to the "models" sheet, recording the
# run prepared for experiments/2025-11-18/0001/run.py.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("models")

$ws.Cells.Item(10, 1).Value = "2025-11-18 15:48:05"
$ws.Cells.Item(10, 2).Value = "models"
$ws.Cells.Item(10, 3).Value = "2025-11-18/a/0001"
$ws.Cells.Item(10, 4).Value = ""
$ws.Cells.Item(10, 5).Value = "models.networks.FCN"
$ws.Cells.Item(10, 6).Value = ""
$ws.Cells.Item(10, 7).Value = "[CallableConfig(path='torch.nn.modules.activation.ReLU', args_cfg=ReLUConfig(inplace=False), kind='class', recovery_mode='call', locked=False, if_recover_while_locked='print')]"
$ws.Cells.Item(10, 8).Value = "[None]"
$ws.Cells.Item(10, 9).Value = "torch.nn.modules.rnn.RNN"
$ws.Cells.Item(10, 10).Value = "embedding_dim___"
$ws.Cells.Item(10, 11).Value = 20
$ws.Cells.Item(10, 12).Value = "tanh"
$ws.Cells.Item(10, 13).Value = "models.networks.FCN"
$ws.Cells.Item(10, 14).Value = "[20, 2]"
$ws.Cells.Item(10, 15).Value = "[CallableConfig(path='torch.nn.modules.linear.Identity', args_cfg=IdentityConfig(), kind='class', recovery_mode='call', locked=False, if_recover_while_locked='print')]"
$ws.Cells.Item(10, 16).Value = "[None]"
